# Add 2022-Q4 data: insert a new quarter sheet (copied from 2022-Q3's
# layout/formatting) with updated figures, and record it in the "总计"
# (totals) summary sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q4" sheet -------------------------------------
# Copy the existing "2022-Q3" sheet to preserve its column layout, borders,
# and header styling; Excel inserts the copy immediately before the source
# sheet, which is exactly the "总计, 2022-Q4, 2022-Q3, 2022-Q2" order we need.
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Update the fund figures on the new sheet (columns D-G for rows 2-4). The
# values are stored as text (matching the source data's inlineStr cells), so
# force text entry and then drop the resulting "quote-prefixed" number
# format back off so no stray cell style lingers.
$q4figures = @{ "D" = "1.06"; "E" = "82.28"; "F" = "2.11"; "G" = "0.0224" }
foreach ($col in $q4figures.Keys) {
    for ($row = 2; $row -le 4; $row++) {
        $cell = $q4.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $q4figures[$col]
        $cell.ClearFormats()
    }
}

# --- 2. Update the "总计" (totals) summary sheet ----------------------------
# Insert the new quarter as row 2 (pushing 2022-Q3 / 2022-Q2 down a row) by
# rewriting the three data rows directly.
$total = $wb.Worksheets.Item("总计")

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.07000000000000001

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.04

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.03

# Give the new A4 "总计" row the same style as the existing A2/A3 index
# cells (centered header-row border style) by copying formats only.
$total.Range("A2").Copy() | Out-Null
$total.Range("A4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Restore the original active sheet ("2022-Q2" was the selected tab before
# this edit) now that the new sheet has been inserted ahead of it.
$wb.Worksheets.Item("2022-Q2").Activate()
